$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at position 20, pushing the totals row (20->21)
#    and the footer row (21->22) down by one.
$ws.Rows(20).Insert()

# 2. Turn the new blank row 20 into a proper data row by cloning the
#    formatting/merges of row 19 (the last data row before insertion).
$ws.Range("A19:N19").Copy($ws.Range("A20:N20"))
$ws.Rows(20).RowHeight = 25.5

# 3. Row 20 now holds a duplicate of row 19's data ("مبرد قدم"), which is
#    exactly what we want there; just fix the serial number.
$ws.Range("A20").Value = 17

# 4. Shift the item data (name / H / L / N columns) of rows 15-19 down by
#    one row, cascading from the bottom up so sources are read before
#    being overwritten.
for ($r = 19; $r -ge 16; $r--) {
    $src = $r - 1
    $ws.Range("B$r").Value = $ws.Range("B$src").Value()
    $ws.Range("H$r").Value = $ws.Range("H$src").Value()
    $ws.Range("L$r").Value = $ws.Range("L$src").Value()
    $ws.Range("N$r").Value = $ws.Range("N$src").Value()
}

# 5. Row 15 becomes the newly added item.
$ws.Range("B15").Value = "PULMICORT 0.25MG/ML 20 NEBULIZER VIAL SUSP."
$ws.Range("H15").Value = "0:11"
$ws.Range("L15").Value = -56.4
$ws.Range("N15").Value = "0:0"

# 6. Update the total in the (now shifted) totals row and normalize its
#    row height to match the rest of the data rows.
$ws.Range("K21").Value = 204.6
$ws.Rows(21).RowHeight = 25.5
